$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a column to plain text (Excel would otherwise coerce
# numeric-looking strings like "241.87" into real numbers).
function Set-TextValue($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextValue "D2" "29.327.50"
$ws.Range("E2").Value = "  +0.60%  "

Set-TextValue "D3" "1.874.00"
$ws.Range("E3").Value = "  +0.70%  "

$ws.Range("E4").Value = "  +0.05%  "

Set-TextValue "D5" "0.7124"
$ws.Range("E5").Value = "  +0.76%  "

Set-TextValue "D6" "241.87"
$ws.Range("E6").Value = "  +0.33%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue "D8" "0.3112"
$ws.Range("E8").Value = "  +0.88%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue "D9" "0.07784"
$ws.Range("E9").Value = "  +1.99%  "

Set-TextValue "D10" "25.11"
$ws.Range("E10").Value = "  +1.79%  "

Set-TextValue "D11" "0.08408"
$ws.Range("E11").Value = "  +0.75%  "

Set-TextValue "D12" "1.873.04"
$ws.Range("E12").Value = "  +0.89%  "

Set-TextValue "D13" "5.244"
$ws.Range("E13").Value = "  +1.23%  "

Set-TextValue "D14" "0.7119"
$ws.Range("E14").Value = "  +0.63%  "

Set-TextValue "D15" "91.15"
$ws.Range("E15").Value = "  +0.14%  "

Set-TextValue "D16" "29.332.71"
$ws.Range("E16").Value = "  +0.47%  "

Set-TextValue "D17" "6.090"
$ws.Range("E17").Value = "  +3.11%  "

Set-TextValue "D18" "0.000008214"
$ws.Range("E18").Value = "  +5.30%  "

Set-TextValue "D19" "240.10"
$ws.Range("E19").Value = "  -1.19%  "

$ws.Range("E20").Value = "  +0.94%  "

Set-TextValue "D21" "2.122.03"
$ws.Range("E21").Value = "  +0.11%  "

$ws.Range("E22").Value = "  -0.06%  "

Set-TextValue "D23" "7.765"
$ws.Range("E23").Value = "  -1.25%  "

Set-TextValue "D24" "1.001"
$ws.Range("E24").Value = "  +0.05%  "

Set-TextValue "D25" "0.1598"
$ws.Range("E25").Value = "  +0.84%  "

Set-TextValue "D26" "162.70"
$ws.Range("E26").Value = "  -0.21%  "

Set-TextValue "D27" "9.036"
$ws.Range("E27").Value = "  +1.31%  "

Set-TextValue "D28" "18.50"
$ws.Range("E28").Value = "  +0.26%  "

Set-TextValue "D29" "1.509"
$ws.Range("E29").Value = "  +0.68%  "

Set-TextValue "D30" "4.419"
$ws.Range("E30").Value = "  +0.50%  "

Set-TextValue "D31" "1.290"
$ws.Range("E31").Value = "  -3.62%  "

Set-TextValue "D32" "4.304"
$ws.Range("E32").Value = "  +2.33%  "

Set-TextValue "D33" "0.05297"
$ws.Range("E33").Value = "  +3.13%  "

Set-TextValue "D34" "1.937"
$ws.Range("E34").Value = "  +1.35%  "

$ws.Range("E35").Value = "  +1.33%  "

Set-TextValue "D36" "0.7449"
$ws.Range("E36").Value = "  -6.52%  "

Set-TextValue "D37" "2.702"
$ws.Range("E37").Value = "  +0.76%  "

Set-TextValue "D38" "0.01868"
$ws.Range("E38").Value = "  +1.44%  "

Set-TextValue "D39" "1.228.19"
$ws.Range("E39").Value = "  +5.19%  "

Set-TextValue "D40" "2.728"
$ws.Range("E40").Value = "  +1.24%  "

Set-TextValue "D41" "6.546"
$ws.Range("E41").Value = "  +5.77%  "

Set-TextValue "D42" "110.80"
$ws.Range("E42").Value = "  +8.72%  "

Set-TextValue "D43" "0.8870"
$ws.Range("E43").Value = "  -0.29%  "

Set-TextValue "D44" "72.67"
$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("E45").Value = "  +0.05%  "

Set-TextValue "D46" "2.019.73"
$ws.Range("E46").Value = "  +0.35%  "

Set-TextValue "D47" "1.802"
$ws.Range("E47").Value = "  +2.01%  "

Set-TextValue "D48" "0.5191"
$ws.Range("E48").Value = "  -0.04%  "

$ws.Range("E49").Value = "  +4.87%  "

Set-TextValue "D50" "9.385"
$ws.Range("E50").Value = "  +0.62%  "

Set-TextValue "D51" "0.4319"
$ws.Range("E51").Value = "  +1.14%  "

